$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.272420883178711
$ws.Range("B1").Value = 2.406774997711182
$ws.Range("C1").Value = 4.486362457275391
$ws.Range("D1").Value = 2.666439294815063
$ws.Range("E1").Value = 1.338737010955811
